$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out old rows 2 and 3
$ws.Range("A2:D3").Clear()

# Set new values for row 1
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = "hola"
$ws.Range("C1").Value = "hola mundo"
$ws.Range("D1").Value = "hola"
$ws.Range("E1").Value = $null
$ws.Range("F1").Value = "hola"

$ws.Range("D6").Select()
